$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.459.46"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.566.56"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D5").Value = "'208.58"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'22.10"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "'0.0866"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.789.78"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "1.568.05"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "'0.519"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "'63.58"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "27.458.80"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "'213.56"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").Value = "0.0₃0690"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "'7.26"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "'2.02"
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").Value = "'153.07"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'0.0470"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").Value = "1.376.10"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").Value = "'0.956"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "'0.537"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "'0.824"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'0.975"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "'1.81"
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("D44").Value = "'64.21"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "'2.16"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D47").Value = "1.702.51"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'85.37"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'0.0959"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("E51").Value = "  -0.58%  "
